# Update the "想去人数" (want-to-go count) figures in column F for rows 2-5
# on both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$newValues = @{
    2 = 8658
    3 = 195
    4 = 402
    5 = 31
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $newValues.Keys) {
        $ws.Range("F$row").Value = $newValues[$row]
    }
}
